# Lesson 6 wordlist: re-sequence the English/Japanese vocabulary rows
# to match the reshuffled lesson ordering described in the commit
# ("updated section placeholder for Lesson 13 (3rd ed.)", etc.).
# Row 1 is the header (English / Japanese) and is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "money"
$ws.Cells.Item(2, 2).Value = "お金|おかね"
$ws.Cells.Item(3, 1).Value = "bath"
$ws.Cells.Item(3, 2).Value = "お風呂|おふろ"
$ws.Cells.Item(4, 1).Value = "kanji; Chinese character"
$ws.Cells.Item(4, 2).Value = "漢字|かんじ"
$ws.Cells.Item(5, 1).Value = "textbook"
$ws.Cells.Item(5, 2).Value = "教科書|きょうかしょ"
$ws.Cells.Item(6, 1).Value = "this week"
$ws.Cells.Item(6, 2).Value = "今週|こんしゅう"
$ws.Cells.Item(7, 1).Value = "CD"
$ws.Cells.Item(7, 2).Value = "シーディー"
$ws.Cells.Item(8, 1).Value = "municipal hospital"
$ws.Cells.Item(8, 2).Value = "市民病院|しみんびょういん"
$ws.Cells.Item(9, 1).Value = "shower"
$ws.Cells.Item(9, 2).Value = "シャワー"
$ws.Cells.Item(10, 1).Value = "next"
$ws.Cells.Item(10, 2).Value = "次|つぎ"
$ws.Cells.Item(11, 1).Value = "electricity"
$ws.Cells.Item(11, 2).Value = "電気|でんき"
$ws.Cells.Item(12, 1).Value = "train"
$ws.Cells.Item(12, 2).Value = "電車|でんしゃ"
$ws.Cells.Item(13, 1).Value = "baggage"
$ws.Cells.Item(13, 2).Value = "荷物|にもつ"
$ws.Cells.Item(14, 1).Value = "personal computer"
$ws.Cells.Item(14, 2).Value = "パソコン"
$ws.Cells.Item(15, 1).Value = "page"
$ws.Cells.Item(15, 2).Value = "ページ"
$ws.Cells.Item(16, 1).Value = "window"
$ws.Cells.Item(16, 2).Value = "窓|まど"
$ws.Cells.Item(17, 1).Value = "night"
$ws.Cells.Item(17, 2).Value = "夜|よる"
$ws.Cells.Item(18, 1).Value = "next week"
$ws.Cells.Item(18, 2).Value = "来週|らいしゅう"
$ws.Cells.Item(19, 1).Value = "next year"
$ws.Cells.Item(19, 2).Value = "来年|らいねん"
$ws.Cells.Item(20, 1).Value = "tough (situation)"
$ws.Cells.Item(20, 2).Value = "大変|たいへん（な）"
$ws.Cells.Item(21, 1).Value = "to play; to spend time pleasantly"
$ws.Cells.Item(21, 2).Value = "遊ぶ|あそぶ"
$ws.Cells.Item(22, 1).Value = "to hurry"
$ws.Cells.Item(22, 2).Value = "急ぐ|いそぐ"
$ws.Cells.Item(23, 1).Value = "to take a bath"
$ws.Cells.Item(23, 2).Value = "お風呂に入る|おふろにはいる"
$ws.Cells.Item(24, 1).Value = "to return (a thing)"
$ws.Cells.Item(24, 2).Value = "返す|かえす"
$ws.Cells.Item(25, 1).Value = "to turn off; to erase"
$ws.Cells.Item(25, 2).Value = "消す|けす"
$ws.Cells.Item(26, 1).Value = "to die"
$ws.Cells.Item(26, 2).Value = "死ぬ|しぬ"
$ws.Cells.Item(27, 1).Value = "to sit down"
$ws.Cells.Item(27, 2).Value = "座る|すわる"
$ws.Cells.Item(28, 1).Value = "to stand up"
$ws.Cells.Item(28, 2).Value = "立つ|たつ"
$ws.Cells.Item(29, 1).Value = "to smoke"
$ws.Cells.Item(29, 2).Value = "たばこを吸う|たばこをすう"
$ws.Cells.Item(30, 1).Value = "to use"
$ws.Cells.Item(30, 2).Value = "使う|つかう"
$ws.Cells.Item(31, 1).Value = "to help"
$ws.Cells.Item(31, 2).Value = "手伝う|てつだう"
$ws.Cells.Item(32, 1).Value = "to enter"
$ws.Cells.Item(32, 2).Value = "入る|はいる"
$ws.Cells.Item(33, 1).Value = "to carry; to hold"
$ws.Cells.Item(33, 2).Value = "持つ|もつ"
$ws.Cells.Item(34, 1).Value = "to be absent (from...); to rest"
$ws.Cells.Item(34, 2).Value = "休む|やすむ"
$ws.Cells.Item(35, 1).Value = "to open (something)"
$ws.Cells.Item(35, 2).Value = "開ける|あける"
$ws.Cells.Item(36, 1).Value = "to teach; to instruct"
$ws.Cells.Item(36, 2).Value = "教える|おしえる"
$ws.Cells.Item(37, 1).Value = "to get off"
$ws.Cells.Item(37, 2).Value = "降りる|おりる"
$ws.Cells.Item(38, 1).Value = "to borrow"
$ws.Cells.Item(38, 2).Value = "借りる|かりる"
$ws.Cells.Item(39, 1).Value = "to close (something)"
$ws.Cells.Item(39, 2).Value = "閉める|しめる"
$ws.Cells.Item(40, 1).Value = "to take a shower"
$ws.Cells.Item(40, 2).Value = "シャワーを浴びる|シャワーをあびる"
$ws.Cells.Item(41, 1).Value = "to turn on"
$ws.Cells.Item(41, 2).Value = "つける"
$ws.Cells.Item(42, 1).Value = "to make a phone call"
$ws.Cells.Item(42, 2).Value = "電話をかける|でんわをかける"
$ws.Cells.Item(43, 1).Value = "to forget; to leave behind"
$ws.Cells.Item(43, 2).Value = "忘れる|わすれる"
$ws.Cells.Item(44, 1).Value = "to bring (a person)"
$ws.Cells.Item(44, 2).Value = "連れてくる|つれてくる"
$ws.Cells.Item(45, 1).Value = "to bring (a thing)"
$ws.Cells.Item(45, 2).Value = "持ってくる|もってくる"
$ws.Cells.Item(46, 1).Value = "later on"
$ws.Cells.Item(46, 2).Value = "後で|あとで"
$ws.Cells.Item(47, 1).Value = "(do something) late"
$ws.Cells.Item(47, 2).Value = "遅く|おそく"
$ws.Cells.Item(48, 1).Value = "...because"
$ws.Cells.Item(48, 2).Value = "～から"
$ws.Cells.Item(49, 1).Value = "That would be fine.; That wouldn't be necessary."
$ws.Cells.Item(49, 2).Value = "結構です|けっこうです"
$ws.Cells.Item(50, 1).Value = "right away"
$ws.Cells.Item(50, 2).Value = "すぐ"
$ws.Cells.Item(51, 1).Value = "Really?"
$ws.Cells.Item(51, 2).Value = "本当ですか|ほんとうですか"
$ws.Cells.Item(52, 1).Value = "slowly; leisurely; unhurriedly"
$ws.Cells.Item(52, 2).Value = "ゆっくり"
$ws.Cells.Item(53, 1).Value = "go straight"
$ws.Cells.Item(53, 2).Value = "まっすぐ行く|まっすぐいき"
$ws.Cells.Item(54, 1).Value = "turn left"
$ws.Cells.Item(54, 2).Value = "左に曲がる|ひだりにまがる"
$ws.Cells.Item(55, 1).Value = "turn right"
$ws.Cells.Item(55, 2).Value = "右に曲がる|みぎにまがる"
$ws.Cells.Item(56, 1).Value = "cross the street"
$ws.Cells.Item(56, 2).Value = "道を渡る|みちをわたる"
$ws.Cells.Item(57, 1).Value = "turn left at the second corner"
$ws.Cells.Item(57, 2).Value = "二つ目の角を左に曲がる|ふたつめのかどをひだりにまがる"
$ws.Cells.Item(58, 1).Value = "turn right at the first traffic light"
$ws.Cells.Item(58, 2).Value = "一つ目の信号を右に曲がる|ひとつめのしんごうをみぎにまがる"
$ws.Cells.Item(59, 1).Value = "left side of the street"
$ws.Cells.Item(59, 2).Value = "道の左側|みちのひだりがわ"
$ws.Cells.Item(60, 1).Value = "right side of the street"
$ws.Cells.Item(60, 2).Value = "道の右側|みちのみぎがわ"
$ws.Cells.Item(61, 1).Value = "north"
$ws.Cells.Item(61, 2).Value = "北|きた"
$ws.Cells.Item(62, 1).Value = "east"
$ws.Cells.Item(62, 2).Value = "東|ひがし"
$ws.Cells.Item(63, 1).Value = "south"
$ws.Cells.Item(63, 2).Value = "南|みなみ"
$ws.Cells.Item(64, 1).Value = "west"
$ws.Cells.Item(64, 2).Value = "西|にし"
$ws.Cells.Item(65, 1).Value = "east"
$ws.Cells.Item(65, 2).Value = "東|ひがし"
$ws.Cells.Item(66, 1).Value = "east exit"
$ws.Cells.Item(66, 2).Value = "東口|ひがしぐち"
$ws.Cells.Item(67, 1).Value = "Tokyo"
$ws.Cells.Item(67, 2).Value = "東京|とうきょう"
$ws.Cells.Item(68, 1).Value = "Kanto area"
$ws.Cells.Item(68, 2).Value = "関東|かんとう"
$ws.Cells.Item(69, 1).Value = "the East"
$ws.Cells.Item(69, 2).Value = "東洋|とうよう"
$ws.Cells.Item(70, 1).Value = "west"
$ws.Cells.Item(70, 2).Value = "西|にし"
$ws.Cells.Item(71, 1).Value = "west exit"
$ws.Cells.Item(71, 2).Value = "西口|にしぐち"
$ws.Cells.Item(72, 1).Value = "northwest"
$ws.Cells.Item(72, 2).Value = "北西|ほくせい"
$ws.Cells.Item(73, 1).Value = "Kansai area"
$ws.Cells.Item(73, 2).Value = "関西|かんさい"
$ws.Cells.Item(74, 1).Value = "the West"
$ws.Cells.Item(74, 2).Value = "西洋|せいよう"
$ws.Cells.Item(75, 1).Value = "south"
$ws.Cells.Item(75, 2).Value = "南|みなみ"
$ws.Cells.Item(76, 1).Value = "south exit"
$ws.Cells.Item(76, 2).Value = "南口|みなみぐち"
$ws.Cells.Item(77, 1).Value = "southeast"
$ws.Cells.Item(77, 2).Value = "南東|なんとう"
$ws.Cells.Item(78, 1).Value = "Antarctica"
$ws.Cells.Item(78, 2).Value = "南極|なんきょく"
$ws.Cells.Item(79, 1).Value = "Southeast Asia"
$ws.Cells.Item(79, 2).Value = "東南アジア|とうなんアジア"
$ws.Cells.Item(80, 1).Value = "north"
$ws.Cells.Item(80, 2).Value = "北|きた"
$ws.Cells.Item(81, 1).Value = "north exit"
$ws.Cells.Item(81, 2).Value = "北口|きたぐち"
$ws.Cells.Item(82, 1).Value = "Tohoku area"
$ws.Cells.Item(82, 2).Value = "東北|とうほく"
$ws.Cells.Item(83, 1).Value = "North Pole"
$ws.Cells.Item(83, 2).Value = "北極|ほっきょく"
$ws.Cells.Item(84, 1).Value = "Hokkaido"
$ws.Cells.Item(84, 2).Value = "北海道|ほっかいどう"
$ws.Cells.Item(85, 1).Value = "north exit"
$ws.Cells.Item(85, 2).Value = "北口|きたぐち"
$ws.Cells.Item(86, 1).Value = "mouth"
$ws.Cells.Item(86, 2).Value = "口|くち"
$ws.Cells.Item(87, 1).Value = "population"
$ws.Cells.Item(87, 2).Value = "人口|じんこう"
$ws.Cells.Item(88, 1).Value = "entrance"
$ws.Cells.Item(88, 2).Value = "入り口／入口|いりぐち"
$ws.Cells.Item(89, 1).Value = "to exit"
$ws.Cells.Item(89, 2).Value = "出る|でる"
$ws.Cells.Item(90, 1).Value = "exit"
$ws.Cells.Item(90, 2).Value = "出口|でぐち"
$ws.Cells.Item(91, 1).Value = "to take something out"
$ws.Cells.Item(91, 2).Value = "出す|だす"
$ws.Cells.Item(92, 1).Value = "attendance"
$ws.Cells.Item(92, 2).Value = "出席|しゅっせき"
$ws.Cells.Item(93, 1).Value = "export"
$ws.Cells.Item(93, 2).Value = "輸出|ゆしゅつ"
$ws.Cells.Item(94, 1).Value = "right"
$ws.Cells.Item(94, 2).Value = "右|みぎ"
$ws.Cells.Item(95, 1).Value = "right turn"
$ws.Cells.Item(95, 2).Value = "右折|うせつ"
$ws.Cells.Item(96, 1).Value = "right and left"
$ws.Cells.Item(96, 2).Value = "左右|さゆう"
$ws.Cells.Item(97, 1).Value = "right hand"
$ws.Cells.Item(97, 2).Value = "右手|みぎて"
$ws.Cells.Item(98, 1).Value = "right side"
$ws.Cells.Item(98, 2).Value = "右側|みぎがわ"
$ws.Cells.Item(99, 1).Value = "left"
$ws.Cells.Item(99, 2).Value = "左|ひだり"
$ws.Cells.Item(100, 1).Value = "left turn"
$ws.Cells.Item(100, 2).Value = "左折|させつ"
$ws.Cells.Item(101, 1).Value = "left hand"
$ws.Cells.Item(101, 2).Value = "左手|ひだりて"
$ws.Cells.Item(102, 1).Value = "left-handed"
$ws.Cells.Item(102, 2).Value = "左利き|ひだりきき"
$ws.Cells.Item(103, 1).Value = "five minutes"
$ws.Cells.Item(103, 2).Value = "五分|ごふん"
$ws.Cells.Item(104, 1).Value = "ten minutes"
$ws.Cells.Item(104, 2).Value = "十分|じゅっぷん／じっぷん"
$ws.Cells.Item(105, 1).Value = "oneself"
$ws.Cells.Item(105, 2).Value = "自分|じぶん"
$ws.Cells.Item(106, 1).Value = "to divide"
$ws.Cells.Item(106, 2).Value = "分ける|わける"
$ws.Cells.Item(107, 1).Value = "teacher"
$ws.Cells.Item(107, 2).Value = "先生|せんせい"
$ws.Cells.Item(108, 1).Value = "last week"
$ws.Cells.Item(108, 2).Value = "先週|せんしゅう"
$ws.Cells.Item(109, 1).Value = "ahead"
$ws.Cells.Item(109, 2).Value = "先に|さきに"
$ws.Cells.Item(110, 1).Value = "last month"
$ws.Cells.Item(110, 2).Value = "先月|せんげつ"
$ws.Cells.Item(111, 1).Value = "senior member"
$ws.Cells.Item(111, 2).Value = "先輩|せんぱい"
$ws.Cells.Item(112, 1).Value = "student"
$ws.Cells.Item(112, 2).Value = "学生|がくせい"
$ws.Cells.Item(113, 1).Value = "to be born"
$ws.Cells.Item(113, 2).Value = "生まれる|うまれる"
$ws.Cells.Item(114, 1).Value = "once in a life time"
$ws.Cells.Item(114, 2).Value = "一生に一度|いっしょうにいちど"
$ws.Cells.Item(115, 1).Value = "college student"
$ws.Cells.Item(115, 2).Value = "大学生|だいがくせい"
$ws.Cells.Item(116, 1).Value = "big"
$ws.Cells.Item(116, 2).Value = "大きい|おおきい"
$ws.Cells.Item(117, 1).Value = "tough"
$ws.Cells.Item(117, 2).Value = "大変な|たいへんな"
$ws.Cells.Item(118, 1).Value = "adult"
$ws.Cells.Item(118, 2).Value = "大人|おとな"
$ws.Cells.Item(119, 1).Value = "embassy"
$ws.Cells.Item(119, 2).Value = "大使館|たいしかん"
$ws.Cells.Item(120, 1).Value = "university"
$ws.Cells.Item(120, 2).Value = "大学|だいがく"
$ws.Cells.Item(121, 1).Value = "student"
$ws.Cells.Item(121, 2).Value = "学生|がくせい"
$ws.Cells.Item(122, 1).Value = "school"
$ws.Cells.Item(122, 2).Value = "学校|がっこう"
$ws.Cells.Item(123, 1).Value = "to study"
$ws.Cells.Item(123, 2).Value = "学ぶ|まなぶ"
$ws.Cells.Item(124, 1).Value = "department; faculty"
$ws.Cells.Item(124, 2).Value = "学部|がくぶ"
$ws.Cells.Item(125, 1).Value = "foreign country"
$ws.Cells.Item(125, 2).Value = "外国|がいこく"
$ws.Cells.Item(126, 1).Value = "foreigner"
$ws.Cells.Item(126, 2).Value = "外国人|がいこくじん"
$ws.Cells.Item(127, 1).Value = "outside"
$ws.Cells.Item(127, 2).Value = "外|そと"
$ws.Cells.Item(128, 1).Value = "overseas"
$ws.Cells.Item(128, 2).Value = "海外|かいがい"
$ws.Cells.Item(129, 1).Value = "China"
$ws.Cells.Item(129, 2).Value = "中国|ちゅうごく"
$ws.Cells.Item(130, 1).Value = "country"
$ws.Cells.Item(130, 2).Value = "国|くに"
$ws.Cells.Item(131, 1).Value = "South Korea"
$ws.Cells.Item(131, 2).Value = "韓国|かんこく"
$ws.Cells.Item(132, 1).Value = "the Diet"
$ws.Cells.Item(132, 2).Value = "国会|こっかい"
